$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1665.3077
$ws.Range("I18").Value = 1745.75
$ws.Range("K18").Value = 1745.75
$ws.Range("M18").Value = -1461.75
$ws.Range("H40").Value = 2144.4443
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2185.7144
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2185.7144
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2535.7144
$ws.Range("H74").Value = 6308.1113
$ws.Range("I74").Value = 5833.3335
$ws.Range("K74").Value = 5833.3335
$ws.Range("M74").Value = -4897.3335
$ws.Range("H77").Value = 6308.1113
$ws.Range("I77").Value = 5833.3335
$ws.Range("K77").Value = 29166.6675
$ws.Range("M77").Value = -24486.6675
$ws.Range("H86").Value = 4623
$ws.Range("I86").Value = 3249.75
$ws.Range("J86").Value = 5233.3335
$ws.Range("K86").Value = 3249.75
$ws.Range("L86").Value = 5233.3335
$ws.Range("M86").Value = -2126.75
$ws.Range("N86").Value = -7479.3335
$ws.Range("H89").Value = 4623
$ws.Range("I89").Value = 3249.75
$ws.Range("J89").Value = 5233.3335
$ws.Range("K89").Value = 16248.75
$ws.Range("L89").Value = 26166.6675
$ws.Range("M89").Value = -10632.75
$ws.Range("N89").Value = -37398.6675
$ws.Range("H92").Value = 504.88235
$ws.Range("I92").Value = 355.2857
$ws.Range("J92").Value = 1203
$ws.Range("K92").Value = 355.2857
$ws.Range("L92").Value = 1203
$ws.Range("M92").Value = 892.7143
$ws.Range("N92").Value = -3699
$ws.Range("H98").Value = 1050.75
$ws.Range("I98").Value = 918
$ws.Range("K98").Value = 918
$ws.Range("M98").Value = 580
$ws.Range("H111").Value = 4005
$ws.Range("J111").Value = 2998
$ws.Range("L111").Value = 8994
$ws.Range("N111").Value = -15128
$ws.Range("H122").Value = 1050.75
$ws.Range("I122").Value = 918
$ws.Range("K122").Value = 2754
$ws.Range("M122").Value = -304
$ws.Range("H132").Value = 1373.2222
$ws.Range("I132").Value = 1373.2222
$ws.Range("K132").Value = 4119.6666
$ws.Range("M132").Value = -1589.6666
$ws.Range("H141").Value = 2161.9
$ws.Range("I141").Value = 1111.125
$ws.Range("J141").Value = 6365
$ws.Range("K141").Value = 3333.375
$ws.Range("L141").Value = 19095
$ws.Range("M141").Value = 1846.625
$ws.Range("N141").Value = -29455

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1788.4375
$ws.Range("I45").Value = 1774.3334
$ws.Range("K45").Value = 1774.3334
$ws.Range("M45").Value = -1397.3334
$ws.Range("H74").Value = 2388.1333
$ws.Range("I74").Value = 1391.3334
$ws.Range("K74").Value = 1391.3334
$ws.Range("M74").Value = -517.3334
$ws.Range("H77").Value = 2388.1333
$ws.Range("I77").Value = 1391.3334
$ws.Range("K77").Value = 6956.666999999999
$ws.Range("M77").Value = -2588.666999999999
$ws.Range("H119").Value = 48797
$ws.Range("J119").Value = 48797
$ws.Range("L119").Value = 48797
$ws.Range("N119").Value = -58473
$ws.Range("H132").Value = 2644.1538
$ws.Range("H135").Value = 41500
$ws.Range("J135").Value = 41500
$ws.Range("L135").Value = 41500
$ws.Range("N135").Value = -51640

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3499
$ws.Range("I105").Value = 2678.9565
$ws.Range("J105").Value = 8214.25
$ws.Range("K105").Value = 2678.9565
$ws.Range("L105").Value = 8214.25
$ws.Range("M105").Value = -931.9564999999998
$ws.Range("N105").Value = -11708.25
$ws.Range("H134").Value = 2470
$ws.Range("I134").Value = 2572.8572
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 7718.571599999999
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -5183.571599999999
$ws.Range("N134").Value = -10320

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 727
$ws.Range("I22").Value = 560.5
$ws.Range("J22").Value = 860.2
$ws.Range("K22").Value = 560.5
$ws.Range("L22").Value = 860.2
$ws.Range("M22").Value = -210.5
$ws.Range("N22").Value = -1560.2
$ws.Range("H99").Value = 16257.95
$ws.Range("I99").Value = 14153.091
$ws.Range("K99").Value = 14153.091
$ws.Range("M99").Value = -12655.091
$ws.Range("H126").Value = 16257.95
$ws.Range("I126").Value = 14153.091
$ws.Range("K126").Value = 42459.273
$ws.Range("M126").Value = -39989.273

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8602.799999999999
$ws.Range("I3").Value = 1004.6667
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 3014.0001
$ws.Range("L3").Value = 60000
$ws.Range("M3").Value = -2902.0001
$ws.Range("N3").Value = -60224
$ws.Range("H7").Value = 12500060
$ws.Range("J7").Value = 68.59999999999999
$ws.Range("L7").Value = 205.8
$ws.Range("N7").Value = -429.8
$ws.Range("H107").Value = 613.67566
$ws.Range("I107").Value = 796
$ws.Range("J107").Value = 603.25714
$ws.Range("K107").Value = 2388
$ws.Range("L107").Value = 1809.77142
$ws.Range("M107").Value = -468
$ws.Range("N107").Value = -5649.77142
$ws.Range("H139").Value = 4919.143
$ws.Range("I139").Value = 2281.3333
$ws.Range("J139").Value = 8436.223
$ws.Range("K139").Value = 6843.999899999999
$ws.Range("L139").Value = 25308.669
$ws.Range("M139").Value = -1703.999899999999
$ws.Range("N139").Value = -35588.669

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 7000
$ws.Range("I9").Value = 15000
$ws.Range("J9").Value = 3000
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = -14830
$ws.Range("N9").Value = -3340
$ws.Range("H70").Value = 5599.8
$ws.Range("J70").Value = 4999.75
$ws.Range("L70").Value = 4999.75
$ws.Range("N70").Value = -5539.75
$ws.Range("H73").Value = 5599.8
$ws.Range("J73").Value = 4999.75
$ws.Range("L73").Value = 4999.75
$ws.Range("N73").Value = -6871.75
$ws.Range("H126").Value = 4942.1113
$ws.Range("I126").Value = 4989.5
$ws.Range("J126").Value = 4928.5713
$ws.Range("K126").Value = 14968.5
$ws.Range("L126").Value = 14785.7139
$ws.Range("M126").Value = -12498.5
$ws.Range("N126").Value = -19725.7139

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 50000204
$ws.Range("I13").Value = 50000204
$ws.Range("K13").Value = 50000204
$ws.Range("M13").Value = -50000064
$ws.Range("H16").Value = 3949
$ws.Range("I16").Value = 3949
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3949
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3779
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 3725
$ws.Range("J46").Value = 4466.6665
$ws.Range("L46").Value = 4466.6665
$ws.Range("N46").Value = -4842.6665
$ws.Range("H61").Value = 3314.7334
$ws.Range("I61").Value = 4203
$ws.Range("J61").Value = 2299.5715
$ws.Range("K61").Value = 4203
$ws.Range("L61").Value = 2299.5715
$ws.Range("M61").Value = -4001
$ws.Range("N61").Value = -2703.5715
$ws.Range("H82").Value = 7471.75
$ws.Range("J82").Value = 6999.6665
$ws.Range("L82").Value = 6999.6665
$ws.Range("N82").Value = -7721.6665
$ws.Range("H85").Value = 7471.75
$ws.Range("J85").Value = 6999.6665
$ws.Range("L85").Value = 6999.6665
$ws.Range("N85").Value = -9495.666499999999
$ws.Range("H113").Value = 3314.7334
$ws.Range("I113").Value = 4203
$ws.Range("J113").Value = 2299.5715
$ws.Range("K113").Value = 4203
$ws.Range("L113").Value = 2299.5715
$ws.Range("M113").Value = -2033
$ws.Range("N113").Value = -6639.5715

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H100").Value = 1668.7273
$ws.Range("I100").Value = 1707.125
$ws.Range("J100").Value = 1566.3334
$ws.Range("K100").Value = 3414.25
$ws.Range("L100").Value = 3132.6668
$ws.Range("M100").Value = -2873.25
$ws.Range("N100").Value = -4214.6668
$ws.Range("H113").Value = 622.61536
$ws.Range("I113").Value = 719.6
$ws.Range("K113").Value = 2158.8
$ws.Range("M113").Value = 11.19999999999982
$ws.Range("H119").Value = 11699
$ws.Range("J119").Value = 11699
$ws.Range("L119").Value = 11699
$ws.Range("N119").Value = -21375
